$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0.1424816079850473
$ws.Cells.Item(2, 4).Value = 0.01197498982678269
$ws.Cells.Item(2, 5).Value = 0.9637622992196953
$ws.Cells.Item(2, 6).Value = 0.318272313789997
$ws.Cells.Item(2, 7).Value = 0.00234217694929143
$ws.Cells.Item(2, 13).Value = 5.46505141147162
$ws.Cells.Item(2, 15).Value = 0.9179566526870815

$ws.Cells.Item(3, 2).Value = 0.1329196949319993
$ws.Cells.Item(3, 4).Value = 0.01043587212457453
$ws.Cells.Item(3, 5).Value = 0.846967803031248
$ws.Cells.Item(3, 6).Value = 0.3139561509140876
$ws.Cells.Item(3, 7).Value = 0.002346686016869877
$ws.Cells.Item(3, 13).Value = 4.77459749794798
$ws.Cells.Item(3, 15).Value = 0.9186121314335765

$ws.Cells.Item(4, 2).Value = 0.1271191138327055
$ws.Cells.Item(4, 4).Value = 0.00948962258665631
$ws.Cells.Item(4, 5).Value = 0.7755721030708003
$ws.Cells.Item(4, 6).Value = 0.3120207158198411
$ws.Cells.Item(4, 7).Value = 0.002349581902428424
$ws.Cells.Item(4, 13).Value = 4.349330696835921
$ws.Cells.Item(4, 15).Value = 0.92148148458881

$ws.Cells.Item(5, 2).Value = 0.1247731336548412
$ws.Cells.Item(5, 4).Value = 0.009103693151232051
$ws.Cells.Item(5, 5).Value = 0.7465529706113472
$ws.Cells.Item(5, 6).Value = 0.3114081520600536
$ws.Cells.Item(5, 7).Value = 0.002350794154623732
$ws.Cells.Item(5, 13).Value = 4.175685366025391
$ws.Cells.Item(5, 15).Value = 0.9232581234990391

$ws.Cells.Item(6, 2).Value = 0.1243846633026209
$ws.Cells.Item(6, 4).Value = 0.009039589826421945
$ws.Cells.Item(6, 5).Value = 0.7417387875898527
$ws.Cells.Item(6, 6).Value = 0.3113169706404477
$ws.Cells.Item(6, 7).Value = 0.002350997394377274
$ws.Cells.Item(6, 13).Value = 4.146830489529577
$ws.Cells.Item(6, 15).Value = 0.9235894381681362

$ws.Cells.Item(7, 2).Value = 0.1270874029223847
$ws.Cells.Item(7, 4).Value = 0.00948441914231779
$ws.Cells.Item(7, 5).Value = 0.7751804419280575
$ws.Cells.Item(7, 6).Value = 0.3120117460874638
$ws.Cells.Item(7, 7).Value = 0.002349598120892267
$ws.Cells.Item(7, 13).Value = 4.34699026974215
$ws.Cells.Item(7, 15).Value = 0.9215030028678086

$ws.Cells.Item(8, 2).Value = 0.139170078435825
$ws.Cells.Item(8, 4).Value = 0.01144454410324158
$ws.Cells.Item(8, 5).Value = 0.9234230072218281
$ws.Cells.Item(8, 6).Value = 0.3166336020988041
$ws.Cells.Item(8, 7).Value = 0.002343705343188862
$ws.Cells.Item(8, 13).Value = 5.227249714439523
$ws.Cells.Item(8, 15).Value = 0.9176628287809763

$ws.Cells.Item(9, 2).Value = 0.1634213682550438
$ws.Cells.Item(9, 4).Value = 0.01527983632681185
$ws.Cells.Item(9, 5).Value = 1.216876692427974
$ws.Cells.Item(9, 6).Value = 0.3315306450376951
$ws.Cells.Item(9, 7).Value = 0.002333152887200454
$ws.Cells.Item(9, 13).Value = 6.94366088030165
$ws.Cells.Item(9, 15).Value = 0.9302937363983119

$ws.Cells.Item(10, 2).Value = 0.1815783244638851
$ws.Cells.Item(10, 4).Value = 0.01809477232907852
$ws.Cells.Item(10, 5).Value = 1.434541918741076
$ws.Cells.Item(10, 6).Value = 0.3462594489764967
$ws.Cells.Item(10, 7).Value = 0.002326001798178814
$ws.Cells.Item(10, 13).Value = 8.200039958094635
$ws.Cells.Item(10, 15).Value = 0.9526842081371001

$ws.Cells.Item(11, 2).Value = 0.1899123572507051
$ws.Cells.Item(11, 4).Value = 0.01937537405569856
$ws.Cells.Item(11, 5).Value = 1.534108347490417
$ws.Cells.Item(11, 6).Value = 0.3538315292143608
$ws.Cells.Item(11, 7).Value = 0.002322877079740379
$ws.Cells.Item(11, 13).Value = 8.770916206805339
$ws.Cells.Item(11, 15).Value = 0.9658967734918065

$ws.Cells.Item(12, 2).Value = 0.1930789124078132
$ws.Cells.Item(12, 4).Value = 0.01986036899103283
$ws.Cells.Item(12, 5).Value = 1.571898849787345
$ws.Cells.Item(12, 6).Value = 0.3568283848104983
$ws.Cells.Item(12, 7).Value = 0.002321712112936933
$ws.Cells.Item(12, 13).Value = 8.987025331649875
$ws.Cells.Item(12, 15).Value = 0.9713502778321299

$ws.Cells.Item(13, 2).Value = 0.1923964657375166
$ws.Cells.Item(13, 4).Value = 0.01975591320163517
$ws.Cells.Item(13, 5).Value = 1.563755994130474
$ws.Cells.Item(13, 6).Value = 0.356177133997349
$ws.Cells.Item(13, 7).Value = 0.00232196219815888
$ws.Cells.Item(13, 13).Value = 8.940485021672714
$ws.Cells.Item(13, 15).Value = 0.970155505533512

$ws.Cells.Item(14, 2).Value = 0.1901726587289261
$ws.Cells.Item(14, 4).Value = 0.01941527340279237
$ws.Cells.Item(14, 5).Value = 1.537215606865459
$ws.Cells.Item(14, 6).Value = 0.3540754637772636
$ws.Cells.Item(14, 7).Value = 0.002322780871312196
$ws.Cells.Item(14, 13).Value = 8.788696871112904
$ws.Cells.Item(14, 15).Value = 0.9663363277464327

$ws.Cells.Item(15, 2).Value = 0.1888118966450065
$ws.Cells.Item(15, 4).Value = 0.01920663089306629
$ws.Cells.Item(15, 5).Value = 1.520970408563386
$ws.Cells.Item(15, 6).Value = 0.3528051140997377
$ws.Cells.Item(15, 7).Value = 0.002323284710058416
$ws.Cells.Item(15, 13).Value = 8.6957140062608
$ws.Cells.Item(15, 15).Value = 0.9640560462531482

$ws.Cells.Item(16, 2).Value = 0.1810351671241079
$ws.Cells.Item(16, 4).Value = 0.01801108745956981
$ws.Cells.Item(16, 5).Value = 1.428046730771342
$ws.Cells.Item(16, 6).Value = 0.3457824889157735
$ws.Cells.Item(16, 7).Value = 0.002326208574936978
$ws.Cells.Item(16, 13).Value = 8.162720727974602
$ws.Cells.Item(16, 15).Value = 0.9518829162831537

$ws.Cells.Item(17, 2).Value = 0.1762833966631803
$ws.Cells.Item(17, 4).Value = 0.01727771149791124
$ws.Cells.Item(17, 5).Value = 1.371187169134004
$ws.Cells.Item(17, 6).Value = 0.3417005391649894
$ws.Cells.Item(17, 7).Value = 0.002328035031505262
$ws.Cells.Item(17, 13).Value = 7.835595490683943
$ws.Cells.Item(17, 15).Value = 0.9452009199867177

$ws.Cells.Item(18, 2).Value = 0.173557304117935
$ws.Cells.Item(18, 4).Value = 0.0168558966666339
$ws.Cells.Item(18, 5).Value = 1.338534258929116
$ws.Cells.Item(18, 6).Value = 0.3394345255493079
$ws.Cells.Item(18, 7).Value = 0.002329097651581692
$ws.Cells.Item(18, 13).Value = 7.647378222436089
$ws.Cells.Item(18, 15).Value = 0.9416416245207984

$ws.Cells.Item(19, 2).Value = 0.1726355000227926
$ws.Cells.Item(19, 4).Value = 0.01671307666192234
$ws.Cells.Item(19, 5).Value = 1.327487110481826
$ws.Cells.Item(19, 6).Value = 0.3386812204071674
$ws.Cells.Item(19, 7).Value = 0.002329459517671642
$ws.Cells.Item(19, 13).Value = 7.583639459917435
$ws.Cells.Item(19, 15).Value = 0.9404848372798256

$ws.Cells.Item(20, 2).Value = 0.1767885066547592
$ws.Cells.Item(20, 4).Value = 0.01735577991538406
$ws.Cells.Item(20, 5).Value = 1.377234612101006
$ws.Cells.Item(20, 6).Value = 0.3421265725708906
$ws.Cells.Item(20, 7).Value = 0.002327839351870309
$ws.Cells.Item(20, 13).Value = 7.870424926547344
$ws.Cells.Item(20, 15).Value = 0.94588272853278

$ws.Cells.Item(21, 2).Value = 0.1908255564278534
$ws.Cells.Item(21, 4).Value = 0.01951532555633406
$ws.Cells.Item(21, 5).Value = 1.545008738945683
$ws.Cells.Item(21, 6).Value = 0.3546892284007868
$ws.Cells.Item(21, 7).Value = 0.002322539912082193
$ws.Cells.Item(21, 13).Value = 8.833282399298469
$ws.Cells.Item(21, 15).Value = 0.9674457746538394

$ws.Cells.Item(22, 2).Value = 0.2000615659746785
$ws.Cells.Item(22, 4).Value = 0.02092708333249504
$ws.Cells.Item(22, 5).Value = 1.655169199612573
$ws.Cells.Item(22, 6).Value = 0.3636562700408916
$ws.Cells.Item(22, 7).Value = 0.002319182990742245
$ws.Cells.Item(22, 13).Value = 9.46217477087805
$ws.Cells.Item(22, 15).Value = 0.9841696429247406

$ws.Cells.Item(23, 2).Value = 0.195126477606081
$ws.Cells.Item(23, 4).Value = 0.0201735502306235
$ws.Cells.Item(23, 5).Value = 1.596325028318375
$ws.Cells.Item(23, 6).Value = 0.3587997883342524
$ws.Cells.Item(23, 7).Value = 0.00232096494542952
$ws.Cells.Item(23, 13).Value = 9.126549879839672
$ws.Cells.Item(23, 15).Value = 0.9749980540788385

$ws.Cells.Item(24, 2).Value = 0.1765601284718485
$ws.Cells.Item(24, 4).Value = 0.01732048577040501
$ws.Cells.Item(24, 5).Value = 1.374500450314628
$ws.Cells.Item(24, 6).Value = 0.3419337116478971
$ws.Cells.Item(24, 7).Value = 0.002327927779420046
$ws.Cells.Item(24, 13).Value = 7.854679003355102
$ws.Cells.Item(24, 15).Value = 0.945573604796607

$ws.Cells.Item(25, 2).Value = 0.156801140802969
$ws.Cells.Item(25, 4).Value = 0.01424297875612979
$ws.Cells.Item(25, 5).Value = 1.137157290901428
$ws.Cells.Item(25, 6).Value = 0.3268527275687632
$ws.Cells.Item(25, 7).Value = 0.002335901203124214
$ws.Cells.Item(25, 13).Value = 6.480247635590501
$ws.Cells.Item(25, 15).Value = 0.9246342974952881
